# Updates the Feuil1 UML-style diagram sheet:
#   - reshuffles class-box labels/text in columns I, K, M
#   - adds a small "Rang" legend in column P
#   - adds new fill colours (box headers vs. box bodies) and a bold header
#     font, mirroring how the author recoloured the class boxes
#   - widens columns K and M to fit the new (longer) text
#   - updates print setup / selection left behind by the edit

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Text content changes
# ---------------------------------------------------------------------

# Column I - "Personnage" box attributes swap + box below becomes the
# "iarmable" interface, and the old Equipement/Armes/Parer list is
# replaced by a new Iequipable/TypeEquipement/Niveau/AugmenterRang list.
$ws.Range("I4").Value = "typeelement"
$ws.Range("I5").Value = "Rang"
$ws.Range("I6").ClearContents()
$ws.Range("I9").ClearContents()
$ws.Range("I10").Value = "iarmable"
$ws.Range("I11").Value = "Attaquer()*"
$ws.Range("I12").Value = "Defendre()*"
$ws.Range("I13").ClearContents()
$ws.Range("I14").ClearContents()
$ws.Range("I15").ClearContents()
$ws.Range("I17").Value = "Iequipable"
$ws.Range("I18").Value = "TypeEquipement"
$ws.Range("I19").Value = "Niveau"
$ws.Range("I20").Value = "AugmenterRang()"

# Column K - "Personnage" class box title + attributes, plus two new
# "Armes"/"Armures" interface boxes further down.
$ws.Range("K2").Value = "Personnage*(Element)"
$ws.Range("K3").Value = "nom"
$ws.Range("K4").Value = "typeelement"
$ws.Range("K5").ClearContents()
$ws.Range("K6").Value = "Liste<Equipement>"
$ws.Range("K9").Value = "Armes: Iarmable, Iequipable"
$ws.Range("K14").Value = "Armures: Iarmable, Iequipable"

# Column M - "Joueur" box title gains its inheritance/interface note;
# the old duplicated Attaquer()*/Defendre()* rows under Monstres go away.
$ws.Range("M2").Value = "Joueur*(Personnage), iarmable"
$ws.Range("M11").ClearContents()
$ws.Range("M12").ClearContents()

# Column P - small new "Rang" enum legend.
$ws.Range("P2").Value = "Rang"
$ws.Range("P3").Value = "Normal"
$ws.Range("P4").Value = "Difficile"
$ws.Range("P5").Value = "Expert"

# ---------------------------------------------------------------------
# 2) Colours / fonts for the class-box headers and bodies
#    (order below matches the order the new fills/fonts/styles were
#    first used so the generated style table lines up with the edit)
# ---------------------------------------------------------------------

# Box headers in column K -> blue fill
$ws.Range("K2,K9,K14").Interior.Color = 13998939

# "Element" header cell in column I -> bold + light grey fill
$ws.Range("I2").Font.Bold = $true
$ws.Range("I2").Interior.Color = 14277081

# Box headers in column M -> gold fill
$ws.Range("M2,M9").Interior.Color = 10086143

# Box bodies in column K (including the two new blank placeholder rows
# in each of the Armes/Armures interface boxes) -> light grey fill
$ws.Range("K3,K4,K10,K11,K15,K16").Interior.Color = 14277081

# Box bodies in column M + the iarmable header in column I -> dark grey fill
$ws.Range("M3,M4,I10").Interior.Color = 10921638

# ---------------------------------------------------------------------
# 3) Column widths (K and M grew to fit the new longer labels)
# ---------------------------------------------------------------------
$ws.Columns("K").ColumnWidth = 20.71
$ws.Columns("M").ColumnWidth = 19.14

# ---------------------------------------------------------------------
# 4) Print setup
# ---------------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# 5) Leave the selection where the author last left it
# ---------------------------------------------------------------------
$ws.Range("M14").Select()
